$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update feature-status marks in column B (GraFx Publisher column)
# Cells that moved from "in progress" (❇️) to "done" (✅)
$ws.Range("B6").Value2  = "✅"
$ws.Range("B14").Value2 = "✅"
$ws.Range("B22").Value2 = "✅"
$ws.Range("B23").Value2 = "✅"
$ws.Range("B57").Value2 = "✅"
$ws.Range("B75").Value2 = "✅"
$ws.Range("B162").Value2 = "✅"
$ws.Range("B163").Value2 = "✅"

# Cells that moved from blank to "in progress" (❇️)
$ws.Range("B62").Value2  = "❇️"
$ws.Range("B84").Value2  = "❇️"
$ws.Range("B85").Value2  = "❇️"
$ws.Range("B86").Value2  = "❇️"
$ws.Range("B89").Value2  = "❇️"
$ws.Range("B90").Value2  = "❇️"
$ws.Range("B91").Value2  = "❇️"
$ws.Range("B94").Value2  = "❇️"
$ws.Range("B111").Value2 = "❇️"
$ws.Range("B112").Value2 = "❇️"
$ws.Range("B113").Value2 = "❇️"
$ws.Range("B116").Value2 = "❇️"
$ws.Range("B118").Value2 = "❇️"
$ws.Range("B135").Value2 = "❇️"
$ws.Range("B161").Value2 = "❇️"

# Legend cell text updated
$ws.Range("B2").Value2 = " ❇️ Being developped"

# Reset the view: scroll back to top-left and select B2
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B2").Select()
